# Disaggregation of commodity Copper
#
# 1. Rename the "Copper ores and concentrates" label to "Copper" (row 7,
#    column C on every year sheet).
# 2. For every year sheet, rotate the historical/avg/max values stored in
#    columns D, E, F (rows 5-8) one step to the right: new D = old F,
#    new E = old D, new F = old E.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Update the commodity label in column C, row 7.
    $ws.Cells.Item(7, 3).Value = "Copper"

    # Rotate D/E/F values for each data row (Neodymium, Dysprosium,
    # Copper, Raw silicon).
    foreach ($r in 5, 6, 7, 8) {
        $dVal = $ws.Cells.Item($r, 4).Value()
        $eVal = $ws.Cells.Item($r, 5).Value()
        $fVal = $ws.Cells.Item($r, 6).Value()

        $ws.Cells.Item($r, 4).Value = $fVal
        $ws.Cells.Item($r, 5).Value = $dVal
        $ws.Cells.Item($r, 6).Value = $eVal
    }
}
